$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts GAME_NAME etc. down to row 3)
$ws.Rows("2:2").Insert()

# Populate the new row with the failure message translation
$ws.Range("A2").Value = "FAIL_MESSAGE"
$ws.Range("B2").Value = "You broke Linky!"
$ws.Range("C2").Value = "Vous avez cassé Linky !"

# Adjust column widths: B narrower, C keeps previous width
$ws.Columns("B").ColumnWidth = 20.64
$ws.Columns("C").ColumnWidth = 22.95

# Leave selection on B3 as in the target state
$ws.Range("B3").Select()
